# Auto-generated edit script applying the diff's numeric cell changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1059.5625
$ws.Range("I19").Value = 773.06665
$ws.Range("J19").Value = 1312.3529
$ws.Range("K19").Value = 773.06665
$ws.Range("L19").Value = 1312.3529
$ws.Range("M19").Value = -598.06665
$ws.Range("N19").Value = -1662.3529
$ws.Range("H88").Value = 1746
$ws.Range("I88").Value = 1758.625
$ws.Range("J88").Value = 1720.75
$ws.Range("K88").Value = 1758.625
$ws.Range("L88").Value = 1720.75
$ws.Range("M88").Value = -1352.625
$ws.Range("N88").Value = -2532.75
$ws.Range("H91").Value = 1746
$ws.Range("I91").Value = 1758.625
$ws.Range("J91").Value = 1720.75
$ws.Range("K91").Value = 1758.625
$ws.Range("L91").Value = 1720.75
$ws.Range("M91").Value = -354.625
$ws.Range("N91").Value = -4528.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1895.6786
$ws.Range("I2").Value = 1763.16
$ws.Range("K2").Value = 1763.16
$ws.Range("M2").Value = -1650.16
$ws.Range("H88").Value = 1525
$ws.Range("J88").Value = 1175.4
$ws.Range("L88").Value = 1175.4
$ws.Range("N88").Value = -1987.4
$ws.Range("H91").Value = 1525
$ws.Range("J91").Value = 1175.4
$ws.Range("L91").Value = 1175.4
$ws.Range("N91").Value = -3983.4
$ws.Range("H116").Value = 1895.6786
$ws.Range("I116").Value = 1763.16
$ws.Range("K116").Value = 1763.16
$ws.Range("M116").Value = 530.8399999999999
$ws.Range("H132").Value = 2831
$ws.Range("I132").Value = 2701.9
$ws.Range("K132").Value = 8105.700000000001
$ws.Range("M132").Value = -5575.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1895.6786
$ws.Range("I3").Value = 1763.16
$ws.Range("K3").Value = 1763.16
$ws.Range("M3").Value = -1649.16
$ws.Range("H86").Value = 3365.1738
$ws.Range("I86").Value = 3122.8438
$ws.Range("J86").Value = 3919.0715
$ws.Range("K86").Value = 3122.8438
$ws.Range("L86").Value = 3919.0715
$ws.Range("M86").Value = -1999.8438
$ws.Range("N86").Value = -6165.0715
$ws.Range("H89").Value = 3365.1738
$ws.Range("I89").Value = 3122.8438
$ws.Range("J89").Value = 3919.0715
$ws.Range("K89").Value = 15614.219
$ws.Range("L89").Value = 19595.3575
$ws.Range("M89").Value = -9998.219000000001
$ws.Range("N89").Value = -30827.3575
$ws.Range("H94").Value = 1589.1666
$ws.Range("I94").Value = 1427.5
$ws.Range("K94").Value = 1427.5
$ws.Range("M94").Value = -976.5
$ws.Range("H99").Value = 3496.6128
$ws.Range("I99").Value = 3456.4167
$ws.Range("J99").Value = 3634.4285
$ws.Range("K99").Value = 3456.4167
$ws.Range("L99").Value = 3634.4285
$ws.Range("M99").Value = -1958.4167
$ws.Range("N99").Value = -6630.4285
$ws.Range("H107").Value = 665.2
$ws.Range("I107").Value = 707
$ws.Range("J107").Value = 289
$ws.Range("K107").Value = 707
$ws.Range("L107").Value = 289
$ws.Range("M107").Value = 1213
$ws.Range("N107").Value = -4129
$ws.Range("H134").Value = 2827.8333
$ws.Range("I134").Value = 2610.3333
$ws.Range("J134").Value = 3262.8333
$ws.Range("K134").Value = 7830.999899999999
$ws.Range("L134").Value = 9788.499899999999
$ws.Range("M134").Value = -5295.999899999999
$ws.Range("N134").Value = -14858.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1303.3334
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1283.3
$ws.Range("I5").Value = 697
$ws.Range("K5").Value = 2091
$ws.Range("M5").Value = -1979
$ws.Range("H6").Value = 14895
$ws.Range("I6").Value = 20174
$ws.Range("K6").Value = 60522
$ws.Range("M6").Value = -60409
$ws.Range("H107").Value = 448.57144
$ws.Range("I107").Value = 454.69232
$ws.Range("K107").Value = 1364.07696
$ws.Range("M107").Value = 555.9230400000001
$ws.Range("H122").Value = 581.7273
$ws.Range("I122").Value = 300.14285
$ws.Range("J122").Value = 1074.5
$ws.Range("K122").Value = 2701.28565
$ws.Range("L122").Value = 9670.5
$ws.Range("M122").Value = -251.2856500000003
$ws.Range("N122").Value = -14570.5
$ws.Range("H131").Value = 9260616
$ws.Range("I131").Value = 55556416
$ws.Range("K131").Value = 166669248
$ws.Range("M131").Value = -166664208
$ws.Range("H135").Value = 1283.3
$ws.Range("I135").Value = 697
$ws.Range("K135").Value = 6273
$ws.Range("M135").Value = -3738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4550473.5
$ws.Range("J7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("N7").Value = -7224
$ws.Range("H8").Value = 4550473.5
$ws.Range("J8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("N8").Value = -7278
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0
$ws.Range("H97").Value = 2375.0952
$ws.Range("I97").Value = 1707.9667
$ws.Range("J97").Value = 4042.9167
$ws.Range("K97").Value = 1707.9667
$ws.Range("L97").Value = 4042.9167
$ws.Range("M97").Value = -1211.9667
$ws.Range("N97").Value = -5034.9167
$ws.Range("H102").Value = 3036.2246
$ws.Range("I102").Value = 3029.4146
$ws.Range("J102").Value = 3071.125
$ws.Range("K102").Value = 3029.4146
$ws.Range("L102").Value = 3071.125
$ws.Range("M102").Value = -1407.4146
$ws.Range("N102").Value = -6315.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 564.82355
$ws.Range("I16").Value = 593.9375
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = 593.9375
$ws.Range("L16").Value = 99
$ws.Range("M16").Value = -423.9375
$ws.Range("N16").Value = -439
$ws.Range("H40").Value = 3790.6843
$ws.Range("I40").Value = 3680
$ws.Range("J40").Value = 4381
$ws.Range("K40").Value = 3680
$ws.Range("L40").Value = 4381
$ws.Range("M40").Value = -3544
$ws.Range("N40").Value = -4653
$ws.Range("H55").Value = 5611.6
$ws.Range("I55").Value = 605.3913
$ws.Range("J55").Value = 22060.572
$ws.Range("K55").Value = 605.3913
$ws.Range("L55").Value = 22060.572
$ws.Range("M55").Value = -432.3913
$ws.Range("N55").Value = -22406.572
$ws.Range("H61").Value = 1170
$ws.Range("I61").Value = 1132.4445
$ws.Range("J61").Value = 1282.6666
$ws.Range("K61").Value = 1132.4445
$ws.Range("L61").Value = 1282.6666
$ws.Range("M61").Value = -930.4445000000001
$ws.Range("N61").Value = -1686.6666
$ws.Range("H100").Value = 2874.9412
$ws.Range("I100").Value = 2874.9412
$ws.Range("K100").Value = 2874.9412
$ws.Range("M100").Value = -2333.9412
$ws.Range("H113").Value = 1170
$ws.Range("I113").Value = 1132.4445
$ws.Range("J113").Value = 1282.6666
$ws.Range("K113").Value = 1132.4445
$ws.Range("L113").Value = 1282.6666
$ws.Range("M113").Value = 1037.5555
$ws.Range("N113").Value = -5622.6666
$ws.Range("H136").Value = 2184.8845
$ws.Range("I136").Value = 1541.0555
$ws.Range("J136").Value = 3633.5
$ws.Range("K136").Value = 4623.166499999999
$ws.Range("L136").Value = 10900.5
$ws.Range("M136").Value = -2073.166499999999
$ws.Range("N136").Value = -16000.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 14798.88
$ws.Range("I14").Value = 7498.5
$ws.Range("K14").Value = 7498.5
$ws.Range("M14").Value = -7330.5
$ws.Range("H81").Value = 2424.9092
$ws.Range("I81").Value = 2523.1
$ws.Range("K81").Value = 5046.2
$ws.Range("M81").Value = -3985.2
$ws.Range("H84").Value = 2424.9092
$ws.Range("I84").Value = 2523.1
$ws.Range("K84").Value = 25231
$ws.Range("M84").Value = -19927
$ws.Range("H100").Value = 828.1
$ws.Range("I100").Value = 935.5
$ws.Range("K100").Value = 1871
$ws.Range("M100").Value = -1330
$ws.Range("H122").Value = 2567.5
$ws.Range("I122").Value = 2675.25
$ws.Range("K122").Value = 8025.75
$ws.Range("M122").Value = -5575.75

Write-Host "Applied all cell updates"